$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column D
$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # -4162 = xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -eq "T/R1") {
        $cell.Value2 = "T"
    } elseif ($val -eq "Students") {
        $cell.Value2 = "SS"
    }
}
